# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the "Price" column (D) as plain text so that values such as "1.003",
# "0.07850" or "27.465.55" are written back exactly as scraped, instead of
# being reinterpreted by Excel as numbers/dates (which would drop trailing
# zeros or otherwise reformat the text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.465.55"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "1.837.68"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -1.08%  "
$ws.Range("D5").Value = "332.47"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").Value = "0.4612"
$ws.Range("E7").Value = "  -2.84%  "
$ws.Range("D8").Value = "0.3824"
$ws.Range("D9").Value = "46.10"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "0.07850"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").Value = "0.9733"
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").Value = "21.11"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").Value = "1.845.34"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "5.874"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "7.008"
$ws.Range("E15").Value = "  -3.26%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "87.67"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "0.06636"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "0.00001028"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").Value = "16.90"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "27.470.07"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "5.330"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "157.00"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("D27").Value = "19.33"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "2.063"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").Value = "5.327"
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").Value = "118.47"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").Value = "0.9529"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "0.09294"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "3.574"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "5.224"
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("D35").Value = "1.318"
$ws.Range("E35").Value = "  -4.33%  "
$ws.Range("D36").Value = "0.05914"
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Value = "0.02191"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").Value = "8.057"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("D40").Value = "0.5810"
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("D41").Value = "0.1834"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").Value = "10.06"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").Value = "1.243"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").Value = "0.5477"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").Value = "11.94"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").Value = "1.863"
$ws.Range("E46").Value = "  -4.01%  "
$ws.Range("D47").Value = "0.06655"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "109.78"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "1.041"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000289"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -1.20%  "
